$d = $word.ActiveDocument

# 1. Document counter: 926 -> 8
$d.Content.Find.Execute("926", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8", 2)

# 2. Megrendelő (customer) name
$d.Content.Find.Execute("Megrendelő: dfgdfgdfgdfg", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Megrendelő: Példa Péter", 2)

# 3. Cím (address)
$d.Content.Find.Execute("Cím: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cím: 9700 Szombathely Nemlétezik utca. 3", 2)

# 4. Elérhetőség (contact / phone)
$d.Content.Find.Execute("Elérhetőség: telefon  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Elérhetőség: telefon  6301234567", 2)

# 5. Megjegyzés (comment)
$d.Content.Find.Execute("Megjegyzés:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Megjegyzés:None", 2)

# 6. Megnevezés (designation)
$d.Content.Find.Execute("Megnevezés: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Megnevezés: sdsd", 2)

# 7. Típus (type)
$d.Content.Find.Execute("Típus: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Típus: None", 2)

# 8. Modell (model)
$d.Content.Find.Execute("Modell: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Modell: None", 2)

# 9. Hibajelenség (fault symptom)
$d.Content.Find.Execute("Hibajelenség: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hibajelenség: None", 2)

# 10. Tartozékok (accessories)
$d.Content.Find.Execute("Tartozékok: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tartozékok: None", 2)

# 11. Szerviz diagnózis (service diagnosis)
$d.Content.Find.Execute("Szerviz diagnózis: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Szerviz diagnózis: None", 2)
